# Generate Report for Handoff
#
# A new handoff round was generated for every file that is currently
# "Ready for handoff" and has not yet been handed back in a given locale
# (still shows the 0001-01-01 00:00:00 placeholder handback date). For
# those rows we now have a fresh handoff .xlf, so:
#   - Priority ("E") becomes "ht" on the locale sheets
#   - Latest Handoff Datetime ("H") is refreshed on the locale sheets
#   - Latest HO Xliff Generate Date ("G") is refreshed on the Overview sheet
#
# Rows already in translation (ba7d7c5a, row 3) or already handed back
# (af3124e6, row 9) are left untouched.

$wb = $excel.ActiveWorkbook

$rows = @(4, 5, 6, 7, 8, 10)

$zhcnDateTime = "2016-10-19 12:11:25"
$dedeDateTime = "2016-10-19 12:11:37"

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = $dedeDateTime
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = $zhcnDateTime
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = $dedeDateTime
}
